# "Charte Graphique" - style edit for the "Combat interface" section.
#
# The "Boutons pokémon" (row 38-39), "Boutons sélection de pokémon"
# (row 45-46) and "Bouton fuite" (row 53) colour references are updated:
#   - Fond (background)  #FF0000 -> #993300 / #0000FF
#   - Bordure (border)   #990000 -> #662200
#
# Also updates the current selection to reflect where the author was
# working (B17, with the view scrolled so row 16 is at the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# Boutons pokémon
$ws.Range("B38").Value = "#993300"
$ws.Range("B39").Value = "#662200"

# Boutons sélection de pokémon
$ws.Range("B45").Value = "#993300"
$ws.Range("B46").Value = "#662200"

# Bouton fuite
$ws.Range("B53").Value = "#0000FF"

# Scroll/selection state left by the author after the edit.
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B17").Select()
